$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scripts")

# Remove row 31 (trailing leftover formula row no longer part of the table)
$ws.Rows.Item(31).Delete()

# Clear the existing AutoFilter criteria / unhide the filtered-out rows
# before writing new data so the engine doesn't mis-autofit hidden rows.
$ws.AutoFilterMode = $false

# Enter the newly-collected run ids for DataShape2 / Engine1 (previously blank)
$ws.Range("A19").Value = 3608
$ws.Range("A20").Value = 7716
$ws.Range("A21").Value = 5576
$ws.Range("A22").Value = 8384
$ws.Range("A23").Value = 9968

# Enter the "MEM!" marker for the DataShape2 / Engine2 rows (previously blank)
$ws.Range("A24:A29").Value = "MEM!"

# Re-sort the whole data block by B (Engine), then C (DataShape), then D (TrainDataSet)
$rng = $ws.Range("A6:J29")
$key1 = $ws.Range("B6:B29")
$key2 = $ws.Range("C6:C29")
$key3 = $ws.Range("D6:D29")
$rng.Sort($key1, 1, $key2, $null, 1, $key3, 1, 1)

# Re-apply the AutoFilter (no active filter criteria this time)
$ws.Range("A5:J29").AutoFilter()

# Update the frozen panes / view to match the new layout (freeze through row29/col D)
$ws.Range("A1").Select()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("E30").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A25:A29").Select()

Write-Output "done"
